$d = $word.ActiveDocument

# Update the date/weekday heading (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2023-05-28 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-29 Monday", 2) | Out-Null

# Update each multiplication-problem cell by position (row-major order),
# since some new values duplicate other old values elsewhere in the table
# (a global Find/Replace would risk double-substitution / ordering bugs).
$table = $d.Tables.Item(1)
$values = @(
    "69×50=", "82×40=", "19×29=", "93×62=", "87×36=",
    "67×81=", "70×86=", "22×14=", "15×52=", "20×31=",
    "12×84=", "97×99=", "29×57=", "78×84=", "70×50=",
    "35×56=", "66×91=", "79×40=", "11×78=", "86×41=",
    "83×44=", "23×82=", "15×84=", "67×78=", "70×70=",
    "56×55=", "87×99=", "53×30=", "100×50=", "46×21=",
    "27×92=", "64×88=", "93×30=", "83×50=", "43×81=",
    "91×31=", "83×24=", "12×20=", "59×41=", "24×45=",
    "51×84=", "36×35=", "73×24=", "62×46=", "78×76=",
    "13×24=", "55×60=", "33×96=", "47×42=", "88×27=",
    "53×69=", "78×37=", "76×61=", "20×10=", "54×46=",
    "96×40=", "21×10=", "15×21=", "33×92=", "24×11=",
    "28×85=", "97×92=", "22×67=", "36×33=", "64×50=",
    "50×37=", "38×25=", "18×46=", "98×54=", "49×37=",
    "97×92=", "71×69=", "85×21=", "38×22=", "69×11=",
    "74×47=", "61×98=", "51×28=", "52×85=", "72×56=",
    "15×33=", "62×94=", "50×18=", "15×80=", "61×53=",
    "93×59=", "70×81=", "72×89=", "20×67=", "15×94=",
    "41×20=", "42×65=", "55×50=", "77×88=", "33×77=",
    "33×92=", "71×36=", "26×37=", "44×62=", "45×86="
)

$numCols = $table.Columns.Count
$numRows = $table.Rows.Count
$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $table.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"
